$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 9).Value = 'b'
$ws.Cells.Item(2, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(12, 9).Value = 'ba'
$ws.Cells.Item(12, 10).Value = 'Appreciation'
$ws.Cells.Item(19, 9).Value = 'sd'
$ws.Cells.Item(19, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(51, 9).Value = 'sd'
$ws.Cells.Item(51, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(53, 9).Value = 'sv'
$ws.Cells.Item(53, 10).Value = 'Statement-opinion'
$ws.Cells.Item(61, 9).Value = 'ba'
$ws.Cells.Item(61, 10).Value = 'Appreciation'
$ws.Cells.Item(64, 9).Value = 'b'
$ws.Cells.Item(64, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(87, 9).Value = 'sv'
$ws.Cells.Item(87, 10).Value = 'Statement-opinion'
$ws.Cells.Item(95, 9).Value = 'aa'
$ws.Cells.Item(95, 10).Value = 'Agree/Accept'
$ws.Cells.Item(96, 9).Value = 'aa'
$ws.Cells.Item(96, 10).Value = 'Agree/Accept'
$ws.Cells.Item(111, 9).Value = 'b'
$ws.Cells.Item(111, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(117, 9).Value = 'b'
$ws.Cells.Item(117, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(121, 9).Value = 'sd'
$ws.Cells.Item(121, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(125, 9).Value = 'sd'
$ws.Cells.Item(125, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(126, 9).Value = 'sd'
$ws.Cells.Item(126, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(156, 9).Value = 'aa'
$ws.Cells.Item(156, 10).Value = 'Agree/Accept'
$ws.Cells.Item(178, 9).Value = 'aa'
$ws.Cells.Item(178, 10).Value = 'Agree/Accept'
$ws.Cells.Item(187, 9).Value = 'qy'
$ws.Cells.Item(187, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(189, 9).Value = 'sv'
$ws.Cells.Item(189, 10).Value = 'Statement-opinion'
$ws.Cells.Item(207, 9).Value = 'sd'
$ws.Cells.Item(207, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(208, 9).Value = 'b'
$ws.Cells.Item(208, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(219, 9).Value = 'ba'
$ws.Cells.Item(219, 10).Value = 'Appreciation'
$ws.Cells.Item(278, 9).Value = 'sv'
$ws.Cells.Item(278, 10).Value = 'Statement-opinion'
$ws.Cells.Item(282, 9).Value = 'qy'
$ws.Cells.Item(282, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(287, 9).Value = 'qy'
$ws.Cells.Item(287, 10).Value = 'Yes-No-Question'
$ws.Cells.Item(290, 9).Value = 'sv'
$ws.Cells.Item(290, 10).Value = 'Statement-opinion'
$ws.Cells.Item(294, 9).Value = 'aa'
$ws.Cells.Item(294, 10).Value = 'Agree/Accept'
$ws.Cells.Item(298, 9).Value = 'aa'
$ws.Cells.Item(298, 10).Value = 'Agree/Accept'
$ws.Cells.Item(299, 9).Value = 'aa'
$ws.Cells.Item(299, 10).Value = 'Agree/Accept'
$ws.Cells.Item(322, 9).Value = 'ba'
$ws.Cells.Item(322, 10).Value = 'Appreciation'
$ws.Cells.Item(326, 9).Value = 'sd'
$ws.Cells.Item(326, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(343, 9).Value = 'ba'
$ws.Cells.Item(343, 10).Value = 'Appreciation'
$ws.Cells.Item(363, 9).Value = 'sv'
$ws.Cells.Item(363, 10).Value = 'Statement-opinion'
$ws.Cells.Item(367, 9).Value = 'aa'
$ws.Cells.Item(367, 10).Value = 'Agree/Accept'
$ws.Cells.Item(369, 9).Value = 'aa'
$ws.Cells.Item(369, 10).Value = 'Agree/Accept'
$ws.Cells.Item(370, 9).Value = 'sd'
$ws.Cells.Item(370, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(380, 9).Value = 'sd'
$ws.Cells.Item(380, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(388, 9).Value = 'sv'
$ws.Cells.Item(388, 10).Value = 'Statement-opinion'
$ws.Cells.Item(392, 9).Value = 'sv'
$ws.Cells.Item(392, 10).Value = 'Statement-opinion'
$ws.Cells.Item(396, 9).Value = 'sv'
$ws.Cells.Item(396, 10).Value = 'Statement-opinion'
$ws.Cells.Item(409, 9).Value = 'aa'
$ws.Cells.Item(409, 10).Value = 'Agree/Accept'
$ws.Cells.Item(415, 9).Value = 'aa'
$ws.Cells.Item(415, 10).Value = 'Agree/Accept'
$ws.Cells.Item(418, 9).Value = 'aa'
$ws.Cells.Item(418, 10).Value = 'Agree/Accept'
$ws.Cells.Item(420, 9).Value = 'aa'
$ws.Cells.Item(420, 10).Value = 'Agree/Accept'
$ws.Cells.Item(424, 9).Value = '%'
$ws.Cells.Item(424, 10).Value = 'Uninterpretable'
$ws.Cells.Item(430, 9).Value = 'sv'
$ws.Cells.Item(430, 10).Value = 'Statement-opinion'
$ws.Cells.Item(433, 9).Value = 'sd'
$ws.Cells.Item(433, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(442, 9).Value = 'sd'
$ws.Cells.Item(442, 10).Value = 'Statement-non-opinion'
